$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Version: 1.8.1 -> 1.8.2
$ws.Range("B3").Value = "1.8.2"

# Experimental: clear the "true" value (B7)
$ws.Range("B7").ClearContents()

# Date: 2023-10-31 -> 2025-11-18
# Written via a text formula + paste-values so Excel keeps the new date
# as plain text (matching the original text-based cell) instead of
# auto-converting it into a date serial number, while preserving the
# cell's existing style/format.
$ws.Range("B8").Formula = "=""2025-11-18"""
$ws.Range("B8").Copy()
$ws.Range("B8").PasteSpecial(-4163)
$excel.CutCopyMode = $false
